$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Restore the test data that was previously removed (revert test data changes)
$ws.Range("C3").Value = "Economic"
$ws.Range("D3").Value = "Economic_radio_button"
$ws.Range("C4").Value = "Quality of Life"
$ws.Range("D4").Value = "Quality of Life_radio_button"
$ws.Range("C5").Value = "Real-world Evidence"
$ws.Range("D5").Value = "Real-world Evidence_radio_button"
